# Weekly fruit/vegetable price update: insert two new weekly observations
# at the top of the data block (rows 600-601), pushing the existing rows
# down by two (to 602-626).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 600, shifting existing rows 600:624 down to 602:626.
$ws.Rows("600:601").Insert()

# --- New row 600 ---
$ws.Range("A600").Value = 3
$ws.Range("B600").Value = "Femacal de La Calera"
$ws.Range("C600").Value = "Coquimbo"
$ws.Range("D600").Value = 44939
$ws.Range("E600").Value = 5
$ws.Range("F600").Value = 100112032
$ws.Range("G600").Value = "Zapallo italiano"
$ws.Range("H600").Value = "Sin especificar"
$ws.Range("I600").Value = "Primera"
$ws.Range("J600").Value = 125
$ws.Range("K600").Value = 3000
$ws.Range("L600").Value = 3300
$ws.Range("M600").Value = 3144
$ws.Range("N600").Value = "$/caja 36 unidades"
$ws.Range("O600").Value = "Provincia de Quillota"
$ws.Range("P600").Value = 87
$ws.Range("Q600").Value = 36
$ws.Range("R600").Value = "Hortaliza"

# --- New row 601 ---
$ws.Range("A601").Value = 3
$ws.Range("B601").Value = "Femacal de La Calera"
$ws.Range("C601").Value = "Coquimbo"
$ws.Range("D601").Value = 44939
$ws.Range("E601").Value = 5
$ws.Range("F601").Value = 100112032
$ws.Range("G601").Value = "Zapallo italiano"
$ws.Range("H601").Value = "Sin especificar"
$ws.Range("I601").Value = "Primera"
$ws.Range("J601").Value = 170
$ws.Range("K601").Value = 6000
$ws.Range("L601").Value = 6300
$ws.Range("M601").Value = 6141
$ws.Range("N601").Value = "$/caja 60 unidades"
$ws.Range("O601").Value = "Provincia de Quillota"
$ws.Range("P601").Value = 102
$ws.Range("Q601").Value = 60
$ws.Range("R601").Value = "Hortaliza"
